$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "44.218.90"
$ws.Cells.Item(2, 5).Value = "  +1.50%  "

$ws.Cells.Item(3, 4).Value = "2.245.20"
$ws.Cells.Item(3, 5).Value = "  +1.03%  "

$ws.Cells.Item(4, 5).Value = "  +0.11%  "

$ws.Cells.Item(5, 4).Value = "'307.06"
$ws.Cells.Item(5, 5).Value = "  -1.45%  "

$ws.Cells.Item(6, 4).Value = "'95.75"
$ws.Cells.Item(6, 5).Value = "  -1.25%  "

$ws.Cells.Item(7, 4).Value = "'0.573"
$ws.Cells.Item(7, 5).Value = "  +1.51%  "

$ws.Cells.Item(8, 5).Value = "  +0.19%  "

$ws.Cells.Item(9, 4).Value = "'0.527"
$ws.Cells.Item(9, 5).Value = "  -0.47%  "

$ws.Cells.Item(10, 4).Value = "'35.12"
$ws.Cells.Item(10, 5).Value = "  -1.20%  "

$ws.Cells.Item(11, 4).Value = "'0.0814"
$ws.Cells.Item(11, 5).Value = "  -0.16%  "

$ws.Cells.Item(12, 4).Value = "'7.25"
$ws.Cells.Item(12, 5).Value = "  -0.80%  "

$ws.Cells.Item(13, 5).Value = "  +0.31%  "

$ws.Cells.Item(14, 4).Value = "2.587.78"
$ws.Cells.Item(14, 5).Value = "  +1.09%  "

$ws.Cells.Item(15, 4).Value = "2.331.62"
$ws.Cells.Item(15, 5).Value = "  +5.21%  "

$ws.Cells.Item(16, 4).Value = "'0.834"
$ws.Cells.Item(16, 5).Value = "  +0.34%  "

$ws.Cells.Item(17, 4).Value = "'13.62"
$ws.Cells.Item(17, 5).Value = "  -2.68%  "

$ws.Cells.Item(18, 4).Value = "44.035.56"
$ws.Cells.Item(18, 5).Value = "  +1.36%  "

$ws.Cells.Item(19, 4).Value = "0.0₃0971"
$ws.Cells.Item(19, 5).Value = "  +1.75%  "

$ws.Cells.Item(20, 4).Value = "'6.41"
$ws.Cells.Item(20, 5).Value = "  +2.64%  "

$ws.Cells.Item(21, 4).Value = "'12.18"
$ws.Cells.Item(21, 5).Value = "  -5.50%  "

$ws.Cells.Item(22, 4).Value = "'65.56"
$ws.Cells.Item(22, 5).Value = "  +0.86%  "

$ws.Cells.Item(23, 4).Value = "'237.03"
$ws.Cells.Item(23, 5).Value = "  +0.98%  "

$ws.Cells.Item(24, 5).Value = "  +0.57%  "

$ws.Cells.Item(25, 5).Value = "  +0.06%  "

$ws.Cells.Item(26, 5).Value = "  +0.17%  "

$ws.Cells.Item(27, 4).Value = "'10.00"
$ws.Cells.Item(27, 5).Value = "  +0.88%  "

$ws.Cells.Item(28, 4).Value = "'38.16"
$ws.Cells.Item(28, 5).Value = "  +5.25%  "

$ws.Cells.Item(29, 4).Value = "'2.20"
$ws.Cells.Item(29, 5).Value = "  +0.28%  "

$ws.Cells.Item(30, 4).Value = "'5.94"
$ws.Cells.Item(30, 5).Value = "  +0.64%  "

$ws.Cells.Item(31, 4).Value = "'20.18"
$ws.Cells.Item(31, 5).Value = "  +2.19%  "

$ws.Cells.Item(32, 4).Value = "'153.04"
$ws.Cells.Item(32, 5).Value = "  -3.66%  "

$ws.Cells.Item(33, 5).Value = "  -2.15%  "

$ws.Cells.Item(34, 5).Value = "  +3.90%  "

$ws.Cells.Item(35, 5).Value = "  -1.03%  "

$ws.Cells.Item(36, 5).Value = "  +3.32%  "

$ws.Cells.Item(37, 5).Value = "  +0.93%  "

$ws.Cells.Item(38, 5).Value = "  -5.51%  "

$ws.Cells.Item(39, 4).Value = "'3.52"
$ws.Cells.Item(39, 5).Value = "  +0.32%  "

$ws.Cells.Item(40, 4).Value = "'14.60"
$ws.Cells.Item(40, 5).Value = "  -4.43%  "

$ws.Cells.Item(41, 4).Value = "'3.85"
$ws.Cells.Item(41, 5).Value = "  -3.15%  "

$ws.Cells.Item(42, 4).Value = "'0.0298"
$ws.Cells.Item(42, 5).Value = "  -1.91%  "

$ws.Cells.Item(43, 5).Value = "  +0.19%  "

$ws.Cells.Item(44, 4).Value = "1.756.70"
$ws.Cells.Item(44, 5).Value = "  +3.86%  "

$ws.Cells.Item(45, 4).Value = "'83.26"
$ws.Cells.Item(45, 5).Value = "  -0.18%  "

$ws.Cells.Item(46, 5).Value = "  +0.43%  "

$ws.Cells.Item(47, 4).Value = "'100.42"
$ws.Cells.Item(47, 5).Value = "  -0.45%  "

$ws.Cells.Item(48, 4).Value = "'4.95"
$ws.Cells.Item(48, 5).Value = "  -2.48%  "

$ws.Cells.Item(49, 4).Value = "'8.16"
$ws.Cells.Item(49, 5).Value = "  +2.04%  "

$ws.Cells.Item(50, 5).Value = "  -2.07%  "

$ws.Cells.Item(51, 4).Value = "'54.93"
$ws.Cells.Item(51, 5).Value = "  -1.30%  "
